# Decrement the "剩余" (remaining) value in column E by 1 for each data
# row (rows 2-99), except row 36 which is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value2 = $current - 1
    }
}
